$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert a new row for "Pepe Leal FC" (id 1326835) above the current row 7
#        (which holds "C R Juvenal" / 1488983), pushing everything below down by one. ---
$ws.Rows(7).Insert()

$ws.Range("A7").Value = "Pepe Leal FC"
$ws.Range("B7").Value = 1326835

# --- 2. Rebuild every hyperlink in column C (rows 2-20) so the hyperlink refs/relationship
#        ids line up correctly with the newly-shifted rows. The engine does not reflow the
#        existing hyperlink anchors when a row is inserted, so remove them all and re-add
#        them in the correct, final row order. ---
$ws.Hyperlinks.Delete()

$teamIds = @(32966, 184499, 186283, 287965, 1273719, 1326835, 1488983, 1747619, 1867254, 2371918, 2916559, 4088673, 14709358, 14933455, 16411206, 19209079, 20651178, 44810918, 47775950)

for ($i = 0; $i -lt $teamIds.Count; $i++) {
    $row = 2 + $i
    $id = $teamIds[$i]
    $cell = $ws.Range("C" + $row)
    $cell.Value = "https://cartola.globo.com/#!/time/" + $id
    $ws.Hyperlinks.Add($cell, "https://cartola.globo.com/", "!/time/" + $id) | Out-Null
}

# Re-apply the original "Hyperlink" cell style to the whole column so every row keeps using
# the same style as before (the Add call above nudges the style onto a freshly duplicated one).
$ws.Range("C2:C20").Style = "Hyperlink"

Write-Host ("Final dimension: " + $ws.UsedRange.Address())
